# Insert a new weekly data row above the existing row 13, pushing all
# subsequent data rows (old rows 13-103) down by one (new rows 14-104).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("13:13").Insert()

$ws.Cells.Item(13,1).Value  = 10
$ws.Cells.Item(13,2).Value  = 'Vega Modelo de Temuco'
$ws.Cells.Item(13,3).Value  = 'La Araucanía'
$ws.Cells.Item(13,4).Value  = 45149
$ws.Cells.Item(13,5).Value  = 9
$ws.Cells.Item(13,6).Value  = 'Fruta'
$ws.Cells.Item(13,7).Value  = 100107
$ws.Cells.Item(13,8).Value  = 'Otros'
$ws.Cells.Item(13,9).Value  = 100107011
$ws.Cells.Item(13,10).Value = 'Tuna'
$ws.Cells.Item(13,11).Value = 'Sin especificar'
$ws.Cells.Item(13,12).Value = 'Primera'
$ws.Cells.Item(13,13).Value = 100
$ws.Cells.Item(13,14).Value = 32000
$ws.Cells.Item(13,15).Value = 32000
$ws.Cells.Item(13,16).Value = 32000
$ws.Cells.Item(13,17).Value = '$/caja 16 kilos'
$ws.Cells.Item(13,18).Value = 'Provincia de Los Andes'
$ws.Cells.Item(13,19).Value = 2000
$ws.Cells.Item(13,20).Value = 16
